# Add a new column S (2022) that mirrors column R (2021), copying the
# formatting from column R and filling in the new year's data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from R4:R34 into S4:S34 first, so the new cells
# inherit the same number formats / borders / fills as their row's R cell.
$ws.Range("R4:R34").Copy()
$ws.Range("S4:S34").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New values for the 2022 column, keyed by row number.
$values = @{
    4  = 2022
    5  = 16.696653653506477
    6  = 17.71894995601205
    7  = 15.612684844888001
    8  = 15.66812062518596
    9  = 16.652881900156387
    10 = 14.667361954014684
    11 = 16.525244796823369
    12 = 19.119250309028729
    13 = 13.749215987119079
    14 = 14.263200620072119
    15 = 16.214093517712168
    16 = 12.189607205170377
    17 = 14.097780631317802
    18 = 16.597474200848456
    19 = 11.56800988291025
    20 = 11.191263248519153
    21 = 11.361761672735106
    22 = 11.015850216858553
    23 = 15.623145704601036
    24 = 18.344423887154832
    25 = 12.847349120106124
    26 = 13.798472231512836
    27 = 15.128863237337196
    28 = 12.394222749619622
    29 = 21.380402934584232
    30 = 19.968977602899539
    31 = 22.891947678227961
    32 = 28.912046224512313
    33 = 31.246721692820181
    34 = 26.427454495987305
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# Move the active selection to T4, matching the post-edit workbook.
$ws.Range("T4").Select()
